# TC03_checkindate.xlsx - "Updated test cases for testcase between 5 to 9"
#
# - Check-in/Check-out dates (G2/H2) updated and forced to Text format so
#   they stay literal strings instead of being re-interpreted as dates.
# - New "Check in Error" column (K) added with the expected validation
#   message for an invalid (checkout-before-checkin) test case.
# - Columns touched by the edit are re-autofit to their content, and the
#   full sheet is left selected (matches a post-edit "select all" pass).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated check-in / check-out date values, stored as literal text.
$ws.Range("G2:H2").NumberFormat = "@"
$ws.Range("G2").Value = "28/05/2016"
$ws.Range("H2").Value = "26/05/2016"

# New column K: error title + validation message for the bad date test case.
$ws.Range("K1").Value = "Check in Error"
$ws.Range("K2").Value = "Check-In Date shall be before than Check-Out Date"

# Re-fit the columns whose content changed/was added (C picks up a width
# too, since it previously had no explicit column width).
$ws.Range("C1").EntireColumn.AutoFit()
$ws.Range("I1").EntireColumn.AutoFit()
$ws.Range("J1").EntireColumn.AutoFit()
$ws.Range("K1").EntireColumn.AutoFit()

# Leave the whole sheet selected.
$ws.Cells.Select()
